$wb = $excel.ActiveWorkbook

# ALC row 2: Mercury Rising
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 320.58334
$ws.Range("I2").Value = 258.81818
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 258.81818
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -145.81818
$ws.Range("N2").Value = -1226

# ALC row 8: On the Drip
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 5522.25
$ws.Range("I8").Value = 5522.25
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 16566.75
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -16427.75
$ws.Range("N8").ClearContents()

# ALC row 64: Forged from the Void
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 914725.6
$ws.Range("I64").Value = 1253497.8
$ws.Range("J64").Value = 11333.333
$ws.Range("K64").Value = 1253497.8
$ws.Range("L64").Value = 11333.333
$ws.Range("M64").Value = -1253249.8
$ws.Range("N64").Value = -11829.333

# ALC row 67: Dodging the Draft (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 914725.6
$ws.Range("I67").Value = 1253497.8
$ws.Range("J67").Value = 11333.333
$ws.Range("K67").Value = 1253497.8
$ws.Range("L67").Value = 11333.333
$ws.Range("M67").Value = -1252639.8
$ws.Range("N67").Value = -13049.333

# ALC row 70: Consecrating Congregation
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2763
$ws.Range("I70").Value = 1950
$ws.Range("J70").Value = 2943.6667
$ws.Range("K70").Value = 5850
$ws.Range("L70").Value = 8831.000100000001
$ws.Range("M70").Value = -5580
$ws.Range("N70").Value = -9371.000100000001

# ALC row 73: Curbing the Contagion (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2763
$ws.Range("I73").Value = 1950
$ws.Range("J73").Value = 2943.6667
$ws.Range("K73").Value = 5850
$ws.Range("L73").Value = 8831.000100000001
$ws.Range("M73").Value = -4914
$ws.Range("N73").Value = -10703.0001

# ALC row 80: Cleansing the Wicked Humours
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 527.6667
$ws.Range("I80").Value = 315.9
$ws.Range("J80").Value = 678.9286
$ws.Range("K80").Value = 947.6999999999999
$ws.Range("L80").Value = 2036.7858
$ws.Range("M80").Value = 50.30000000000007
$ws.Range("N80").Value = -4032.7858

# ALC row 83: Washing Away the Sins (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 527.6667
$ws.Range("I83").Value = 315.9
$ws.Range("J83").Value = 678.9286
$ws.Range("K83").Value = 2843.1
$ws.Range("L83").Value = 6110.3574
$ws.Range("M83").Value = 2148.9
$ws.Range("N83").Value = -16094.3574

# ALC row 124: Luncheon Bound
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H124").Value = 34660
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 34660
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 34660
$ws.Range("N124").Value = -44480

# ALC row 125: Body over Mind
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 22422742
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 22422742
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 201804678
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -201809598

# ALC row 137: Cutting Edge of Culinary Quality
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 37038750
$ws.Range("I137").Value = 55556732
$ws.Range("J137").Value = 2790.5557
$ws.Range("K137").Value = 166670196
$ws.Range("L137").Value = 8371.667099999999
$ws.Range("M137").Value = -166667646
$ws.Range("N137").Value = -13471.6671

# ARM row 11: Rodents of Unusual Size
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 2250
$ws.Range("I11").Value = 1500
$ws.Range("J11").Value = 3000
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -1356
$ws.Range("N11").Value = -3288

# ARM row 122: Haste for High Durium
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2837.6667
$ws.Range("I122").Value = 2670.6667
$ws.Range("J122").Value = 3004.6667
$ws.Range("K122").Value = 8012.000100000001
$ws.Range("L122").Value = 9014.000100000001
$ws.Range("M122").Value = -5562.000100000001
$ws.Range("N122").Value = -13914.0001

# BSM row 86: Through Thick and Thin
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 27126.75
$ws.Range("I86").Value = 2833.3333
$ws.Range("J86").Value = 100007
$ws.Range("K86").Value = 2833.3333
$ws.Range("L86").Value = 100007
$ws.Range("M86").Value = -1710.3333
$ws.Range("N86").Value = -102253

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 27126.75
$ws.Range("I89").Value = 2833.3333
$ws.Range("J89").Value = 100007
$ws.Range("K89").Value = 14166.6665
$ws.Range("L89").Value = 500035
$ws.Range("M89").Value = -8550.666499999999
$ws.Range("N89").Value = -511267

# CRP row 11: Leaving without Leave
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 25050
$ws.Range("I11").Value = 25050
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 25050
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -24910

# CRP row 62: Splinter in the Sewers
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 20128.785
$ws.Range("I62").Value = 68000
$ws.Range("J62").Value = 7073
$ws.Range("K62").Value = 68000
$ws.Range("L62").Value = 7073
$ws.Range("M62").Value = -67376
$ws.Range("N62").Value = -8321

# CRP row 65: The Lumber of Their Discontent (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 20128.785
$ws.Range("I65").Value = 68000
$ws.Range("J65").Value = 7073
$ws.Range("K65").Value = 340000
$ws.Range("L65").Value = 35365
$ws.Range("M65").Value = -336880
$ws.Range("N65").Value = -41605

# GSM row 7: Water of Life
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2000
$ws.Range("N7").Value = -2224

# GSM row 8: Gods of Small Things
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 2000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2000
$ws.Range("N8").Value = -2278

# GSM row 58: The Big Red
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 5000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 5000
$ws.Range("N58").Value = -5554

# GSM row 80: Needs More Prayerbell
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2799.4
$ws.Range("I80").Value = 2699.5
$ws.Range("J80").Value = 2866
$ws.Range("K80").Value = 2699.5
$ws.Range("L80").Value = 2866
$ws.Range("M80").Value = -1701.5
$ws.Range("N80").Value = -4862

# GSM row 83: With a Noise That Reaches Heaven (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2799.4
$ws.Range("I83").Value = 2699.5
$ws.Range("J83").Value = 2866
$ws.Range("K83").Value = 13497.5
$ws.Range("L83").Value = 14330
$ws.Range("M83").Value = -8505.5
$ws.Range("N83").Value = -24314

# GSM row 97: If I'd a Koppranickel for Every Time...
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1262
$ws.Range("I97").Value = 1327.5
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 1327.5
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -831.5
$ws.Range("N97").Value = -1992

# GSM row 122: Awarding Academic Excellence
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 742227.8
$ws.Range("I122").Value = 856110.4399999999
$ws.Range("J122").Value = 1990.5
$ws.Range("K122").Value = 2568331.32
$ws.Range("L122").Value = 5971.5
$ws.Range("M122").Value = -2565881.32
$ws.Range("N122").Value = -10871.5

# GSM row 126: Gold Rush Order
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2092.3333
$ws.Range("I126").Value = 1691.7693
$ws.Range("J126").Value = 2464.2856
$ws.Range("K126").Value = 5075.3079
$ws.Range("L126").Value = 7392.8568
$ws.Range("M126").Value = -2605.3079
$ws.Range("N126").Value = -12332.8568

# LTW row 7: Tan Before the Ban
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3094.6667
$ws.Range("I7").Value = 1940.8
$ws.Range("J7").Value = 3538.4614
$ws.Range("K7").Value = 1940.8
$ws.Range("L7").Value = 3538.4614
$ws.Range("M7").Value = -1828.8
$ws.Range("N7").Value = -3762.4614

# LTW row 40: Best Served Toad
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2800.158
$ws.Range("I40").Value = 747.6667
$ws.Range("J40").Value = 3185
$ws.Range("K40").Value = 747.6667
$ws.Range("L40").Value = 3185
$ws.Range("M40").Value = -611.6667
$ws.Range("N40").Value = -3457

# LTW row 82: Trainin' the Neck
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 496.4
$ws.Range("I82").Value = 534
$ws.Range("J82").Value = 440
$ws.Range("K82").Value = 534
$ws.Range("L82").Value = 440
$ws.Range("M82").Value = -173
$ws.Range("N82").Value = -1162

# LTW row 85: Training Is Only Skintight (L)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 496.4
$ws.Range("I85").Value = 534
$ws.Range("J85").Value = 440
$ws.Range("K85").Value = 534
$ws.Range("L85").Value = 440
$ws.Range("M85").Value = 714
$ws.Range("N85").Value = -2936

# LTW row 126: Battered Books
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3094.6667
$ws.Range("I126").Value = 1940.8
$ws.Range("J126").Value = 3538.4614
$ws.Range("K126").Value = 5822.4
$ws.Range("L126").Value = 10615.3842
$ws.Range("M126").Value = -3352.4
$ws.Range("N126").Value = -15555.3842

# WVR row 123: Helping Handwear
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 24217.39
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 24217.39
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 24217.39
$ws.Range("N123").Value = -34017.39
